$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; existing rows 16-28 shift down to 17-29
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new weekly price record
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 45126
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = "Frutos de pepita"
$ws.Range("I16").Value = 100104003
$ws.Range("J16").Value = "Membrillo"
$ws.Range("K16").Value = "Champion"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("Q16").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 806
$ws.Range("T16").Value = 18
